$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Six paragraphs in the patient-info block: tab stop 1985 -> 1418 twips,
#    ind right=567 -> left=-142 right=424.  Also drop the stray _GoBack
#    bookmark that sat in the "Nama Pasien" paragraph.
# ---------------------------------------------------------------------------

# Para 10: "Dokter Pengirim"
$d.Paragraphs.Item(10).Range.InsertXML(@"
<w:p $wNs><w:pPr><w:tabs><w:tab w:val="left" w:pos="1418"/></w:tabs><w:spacing w:after="0"/><w:ind w:left="-142" w:right="424"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00B2745C"><w:t xml:space="preserve">Dokter Pengirim </w:t></w:r><w:r><w:tab/></w:r><w:r w:rsidRPr="00B2745C"><w:t>:</w:t></w:r><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> `${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>dr_pengirim</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p>
"@) | Out-Null

# Para 11: "No  RM"
$d.Paragraphs.Item(11).Range.InsertXML(@"
<w:p $wNs><w:pPr><w:tabs><w:tab w:val="left" w:pos="1418"/></w:tabs><w:spacing w:after="0"/><w:ind w:left="-142" w:right="424"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:t>No  RM</w:t></w:r><w:r><w:tab/></w:r><w:r w:rsidRPr="00B2745C"><w:t>:</w:t></w:r><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> `${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>rm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p>
"@) | Out-Null

# Para 12: "Nama Pasien" (bookmark removed)
$d.Paragraphs.Item(12).Range.InsertXML(@"
<w:p $wNs><w:pPr><w:tabs><w:tab w:val="left" w:pos="1418"/></w:tabs><w:spacing w:after="0"/><w:ind w:left="-142" w:right="424"/></w:pPr><w:r w:rsidRPr="00B2745C"><w:t>Nama Pasien</w:t></w:r><w:r w:rsidRPr="00B2745C"><w:tab/><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>`${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>nama</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p>
"@) | Out-Null

# Para 13: "Jk / Umur"
$d.Paragraphs.Item(13).Range.InsertXML(@"
<w:p $wNs><w:pPr><w:tabs><w:tab w:val="left" w:pos="1418"/></w:tabs><w:spacing w:after="0"/><w:ind w:left="-142" w:right="424"/></w:pPr><w:r><w:t>Jk / Umur</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>`${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>jns_kelamin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>} / `${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>umur</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p>
"@) | Out-Null

# Para 14: "Alamat"
$d.Paragraphs.Item(14).Range.InsertXML(@"
<w:p $wNs><w:pPr><w:tabs><w:tab w:val="left" w:pos="1418"/></w:tabs><w:spacing w:after="0"/><w:ind w:left="-142" w:right="424"/></w:pPr><w:r w:rsidRPr="00B2745C"><w:t xml:space="preserve">Alamat    </w:t></w:r><w:r w:rsidRPr="00B2745C"><w:tab/><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>`${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>alamat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p>
"@) | Out-Null

# Para 15: "Tanggal/ Pukul" - tab+ind change, cols space 708 -> 283, merge the
# tab-run and colon-run into a single run.
$d.Paragraphs.Item(15).Range.InsertXML(@"
<w:p $wNs><w:pPr><w:tabs><w:tab w:val="left" w:pos="1418"/></w:tabs><w:spacing w:after="0"/><w:ind w:left="-142" w:right="424"/><w:rPr><w:lang w:val="en-US"/></w:rPr><w:sectPr w:rsidR="00064DEC" w:rsidRPr="00402228" w:rsidSect="00064DEC"><w:type w:val="continuous"/><w:pgSz w:w="12191" w:h="18711"/><w:pgMar w:top="1134" w:right="567" w:bottom="1440" w:left="567" w:header="284" w:footer="709" w:gutter="0"/><w:cols w:num="2" w:space="283"/><w:docGrid w:linePitch="360"/></w:sectPr></w:pPr><w:r><w:t>Tanggal/ Pukul</w:t></w:r><w:r><w:tab/><w:t>:</w:t></w:r><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> `${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>tgl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00402228"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">} </w:t></w:r></w:p>
"@) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Pemeriksa" signature block: drop trailing tab run, add a new PEMERIKSA
#    column (3 blank lines + a centred ${pemeriksa} line), and move the
#    _GoBack bookmark to a fresh trailing paragraph.
# ---------------------------------------------------------------------------

$pemeriksaPara = $d.Paragraphs.Item(39)
$pemeriksaPara.Range.InsertXML(@"
<w:p $wNs><w:pPr><w:spacing w:after="0"/><w:ind w:left="567" w:right="567" w:firstLine="709"/></w:pPr><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">  </w:t></w:r><w:r w:rsidRPr="00B2745C"><w:t>Pemeriksa</w:t></w:r></w:p>
"@)

$lastPara = $d.Paragraphs.Item(40)
$lastPara.Range.InsertXML(@"
<w:p $wNs><w:pPr><w:spacing w:after="0"/><w:ind w:left="567" w:right="567" w:firstLine="709"/></w:pPr></w:p><w:p $wNs><w:pPr><w:spacing w:after="0"/><w:ind w:left="567" w:right="567" w:firstLine="709"/></w:pPr></w:p><w:p $wNs><w:pPr><w:spacing w:after="0"/><w:ind w:left="567" w:right="567" w:firstLine="709"/></w:pPr></w:p><w:p $wNs><w:pPr><w:spacing w:after="0"/><w:ind w:left="6379" w:right="567"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>`${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>pemeriksa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p $wNs><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@)
